# DocentesImportacionMasiva.xlsx — "bug de tema materia vacio"
#
# The "legajo" column (B) was stored as a plain number (90000, 90001, ...).
# Downstream code that expects the provider/professor code as text (e.g.
# matched/concatenated with a "PR" prefix elsewhere) was seeing it as a
# number, producing an empty "materia" (subject) lookup. The fix: re-enter
# each legajo value as literal text "PR" + <original number>, which forces
# Excel to store it as a shared string instead of a numeric value.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")
$ws.Activate()

# Rows 2..62 hold one record each; legajo numbers run 90000..90060 in order.
$firstRow = 2
$lastRow = 62
$startLegajo = 90000

for ($row = $firstRow; $row -le $lastRow; $row++) {
    $legajo = $startLegajo + ($row - $firstRow)
    $ws.Cells.Item($row, 2).Value = "PR$legajo"
}

# Leave the sheet scrolled/selected where the user was last working when the
# fix was verified (around the bottom of the imported range).
$ws.Application.ActiveWindow.ScrollRow = 41
[void]$ws.Range("A54").Select()
